$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.442.73"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.329.43"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "3.328.52"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "3.911.24"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.70"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.62%  "
$ws.Range("D16").Value = "65.475.95"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000168"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "3.332.42"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.519"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.87"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.80"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "2.716.77"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.22"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "335.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0666"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.73"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.92"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.966"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.57%  "
